$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 85
$ws.Range("I5").Value = 91.666664
$ws.Range("J5").Value = 75
$ws.Range("K5").Value = 91.666664
$ws.Range("L5").Value = 75
$ws.Range("M5").Value = 23.333336
$ws.Range("N5").Value = -305
# Row 40
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 5000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 5000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -4825
$ws.Range("N40").Value = $null
# Row 51
$ws.Range("H51").Value = 4105.4814
$ws.Range("I51").Value = 4350.1
$ws.Range("J51").Value = 3961.5881
$ws.Range("K51").Value = 4350.1
$ws.Range("L51").Value = 3961.5881
$ws.Range("M51").Value = -3866.1
$ws.Range("N51").Value = -4929.5881
# Row 64
$ws.Range("H64").Value = 3391.5283
$ws.Range("I64").Value = 3223.75
$ws.Range("K64").Value = 3223.75
$ws.Range("M64").Value = -2975.75
# Row 67
$ws.Range("H67").Value = 3391.5283
$ws.Range("I67").Value = 3223.75
$ws.Range("K67").Value = 3223.75
$ws.Range("M67").Value = -2365.75
# Row 80
$ws.Range("H80").Value = 35732016
$ws.Range("I80").Value = 26233.385
$ws.Range("J80").Value = 66677028
$ws.Range("K80").Value = 78700.155
$ws.Range("L80").Value = 200031084
$ws.Range("M80").Value = -77702.155
$ws.Range("N80").Value = -200033080
# Row 83
$ws.Range("H83").Value = 35732016
$ws.Range("I83").Value = 26233.385
$ws.Range("J83").Value = 66677028
$ws.Range("K83").Value = 236100.465
$ws.Range("L83").Value = 600093252
$ws.Range("M83").Value = -231108.465
$ws.Range("N83").Value = -600103236
# Row 113
$ws.Range("H113").Value = 3397.8
$ws.Range("I113").Value = 2710.85
$ws.Range("J113").Value = 4313.7334
$ws.Range("K113").Value = 2710.85
$ws.Range("L113").Value = 4313.7334
$ws.Range("M113").Value = 543.1500000000001
$ws.Range("N113").Value = -10821.7334
# Row 129
$ws.Range("H129").Value = 1031.9286
$ws.Range("J129").Value = 1137.909
$ws.Range("L129").Value = 3413.727
$ws.Range("N129").Value = -13413.727
# Row 137
$ws.Range("H137").Value = 29038.236
$ws.Range("I137").Value = 70967.53
$ws.Range("J137").Value = 1693.0435
$ws.Range("K137").Value = 212902.59
$ws.Range("L137").Value = 5079.1305
$ws.Range("M137").Value = -210352.59
$ws.Range("N137").Value = -10179.1305

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2522.06
$ws.Range("I32").Value = 2385.3264
$ws.Range("J32").Value = 5120
$ws.Range("K32").Value = 2385.3264
$ws.Range("L32").Value = 5120
$ws.Range("M32").Value = -2098.3264
$ws.Range("N32").Value = -5694
# Row 39
$ws.Range("H39").Value = 1016
$ws.Range("I39").Value = 1016
$ws.Range("K39").Value = 1016
$ws.Range("M39").Value = -496
# Row 52
$ws.Range("H52").Value = 34500
$ws.Range("J52").Value = 34500
$ws.Range("L52").Value = 34500
$ws.Range("N52").Value = -35136

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 901.53705
$ws.Range("I94").Value = 768.9524
$ws.Range("J94").Value = 1365.5834
$ws.Range("K94").Value = 768.9524
$ws.Range("L94").Value = 1365.5834
$ws.Range("M94").Value = -317.9524
$ws.Range("N94").Value = -2267.5834

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 629000
$ws.Range("I4").Value = 80000
$ws.Range("J4").Value = 766250
$ws.Range("K4").Value = 80000
$ws.Range("L4").Value = 766250
$ws.Range("M4").Value = -79888
$ws.Range("N4").Value = -766474
# Row 38
$ws.Range("H38").Value = 7037.6665
$ws.Range("J38").Value = 7865.2
$ws.Range("L38").Value = 7865.2
$ws.Range("N38").Value = -8619.200000000001
# Row 46
$ws.Range("H46").Value = 7037.6665
$ws.Range("J46").Value = 7865.2
$ws.Range("L46").Value = 7865.2
$ws.Range("N46").Value = -8287.200000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 2249.5
$ws.Range("I4").Value = 2249.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 6748.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -6636.5
$ws.Range("N4").Value = $null

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 5800
$ws.Range("I5").Value = 200
$ws.Range("J5").Value = 7666.6665
$ws.Range("K5").Value = 200
$ws.Range("L5").Value = 7666.6665
$ws.Range("M5").Value = -88
$ws.Range("N5").Value = -7890.6665
# Row 15
$ws.Range("H15").Value = 6375
$ws.Range("J15").Value = 6375
$ws.Range("L15").Value = 6375
$ws.Range("N15").Value = -6951
# Row 81
$ws.Range("H81").Value = 6375
$ws.Range("J81").Value = 6375
$ws.Range("L81").Value = 6375
$ws.Range("N81").Value = -8371
# Row 84
$ws.Range("H84").Value = 6375
$ws.Range("J84").Value = 6375
$ws.Range("L84").Value = 19125
$ws.Range("N84").Value = -29109

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 90916260
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 90916260
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 90916260
$ws.Range("M2").Value = $null
$ws.Range("N2").Value = -90916484
# Row 16
$ws.Range("H16").Value = 3391.353
$ws.Range("I16").Value = 2829.25
$ws.Range("J16").Value = 4740.4
$ws.Range("K16").Value = 2829.25
$ws.Range("L16").Value = 4740.4
$ws.Range("M16").Value = -2659.25
$ws.Range("N16").Value = -5080.4
# Row 61
$ws.Range("H61").Value = 2813.5
$ws.Range("I61").Value = 2751.3333
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 2751.3333
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -2549.3333
$ws.Range("N61").Value = -3404
# Row 100
$ws.Range("H100").Value = 25644756
$ws.Range("I100").Value = 4176.24
$ws.Range("J100").Value = 71431500
$ws.Range("K100").Value = 4176.24
$ws.Range("L100").Value = 71431500
$ws.Range("M100").Value = -3635.24
$ws.Range("N100").Value = -71432582
# Row 113
$ws.Range("H113").Value = 2813.5
$ws.Range("I113").Value = 2751.3333
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2751.3333
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -581.3332999999998
$ws.Range("N113").Value = -7340
# Row 136
$ws.Range("H136").Value = 5452.0835
$ws.Range("I136").Value = 3105.5881
$ws.Range("J136").Value = 7551.579
$ws.Range("K136").Value = 9316.764299999999
$ws.Range("L136").Value = 22654.737
$ws.Range("M136").Value = -6766.764299999999
$ws.Range("N136").Value = -27754.737

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 7986.6665
$ws.Range("I2").Value = 7986.6665
$ws.Range("K2").Value = 7986.6665
$ws.Range("M2").Value = -7874.6665
